# Applies the documented edits to analysis_report.docx.
# Uses Paragraph.Range.Text assignment (rather than Find.Execute) so that the
# xml:space="preserve" attribute on the run's <w:t> is retained, matching how
# Word round-trips single-run paragraphs. Range.Text includes a trailing
# carriage return (paragraph mark) which must be stripped before building the
# replacement string, otherwise assigning it back inserts a spurious empty
# paragraph.

$d = $word.ActiveDocument
$CR = [char]13

$replacements = @(
    @{ match = "First I loaded the data using R language, by exploring the data, there are 668 rows and 125 columns."; new = "First, I loaded the data using R language and did some data exploratory analysis. There are 668 rows and 125 columns." },
    @{ match = "From instruction on powerpoint, we need to exclude the international student, so I filtered and excluded them from the dataframe."; new = "From instruction on power point, we need to exclude the international student, so I filtered and excluded them from the data frame." },
    @{ match = "First, I did analysis and find out if there are cases that are completed empty (full of NAs), these records should be excluded for the further analysis."; new = "First, I did analysis and find out if there are cases that are completed empty (full of NAs)in the entire row, these records should be excluded for the further analysis. There are 0 rows has been removed for this step." },
    @{ match = "Second, find out how many NAs in each row and each column, if the NA rate is higher than a cut out rate, for example, 50%, then the row and columns should be excluded."; new = "Second, find out how many NAs in each column, if the NA rate is higher than a cut out rate (eg. >50%) in the column, that columns will be excluded. There are 25 column has been removed for this step." },
    @{ match = "For rows and columns that has lower NA rate, we need to remove NAs and replaced them with a value. Currently, we need imputation. Imputation is especially important in advanced data analysis. There are lots of methods of data imputation, for this analysis, I used averaged imputation. In this way, missing values are taken care of."; new = "For the rest of table still has some NAs, I did the imputation to replace NAs into a value. Currently, i imputed using the mean of the column. Imputation is especially important in data analysis and there are lots of methods for data imputation, however, for this analysis I used averaged imputation. In this way, missing values are taken care of." },
    @{ match = "When compare if there is a significant difference between URG students’ group and non-URG students’ group for the categories in question faculty mentoring and advising, I used pairwise t-test to calculate p value. In those categories: selection of a dissertation topic, your dissertation research, writing and revising your dissertation, academic career option, nonacademic career option, search for employment or training, I found that there is no significant statistical difference between URG and Non-URG students as p value is larger than 0.05."; new = "When compare if there is a significant difference between URG students’ group and non-URG students’ group for the categories, I used pairwise t-test to calculate p value. If P value > 0.05, there is no significant statistical difference between URG and Non-URG students, otherwise p value < 0.05, significant difference between these two groups." },
    @{ match = "it can be a rebalance way and can make survey more accurately to reflect population."; new = "it can be a re-balance way and can make survey more accurately to reflect population." }
)

$applied = @{}

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $text = $p.Range.Text.TrimEnd($CR)
    for ($k = 0; $k -lt $replacements.Count; $k++) {
        $r = $replacements[$k]
        if ($text -eq $r.match) {
            $p.Range.Text = $r.new
            $applied[$k] = $true
            break
        }
        elseif ($text.Contains($r.match)) {
            $p.Range.Text = $text.Replace($r.match, $r.new)
            $applied[$k] = $true
            break
        }
    }
}

for ($k = 0; $k -lt $replacements.Count; $k++) {
    if (-not $applied[$k]) {
        Write-Output "NOT FOUND: $($replacements[$k].match)"
    }
}

# Merge the "I also convert likert..." paragraph with the "In the category..."
# paragraph into a single rewritten paragraph, deleting the redundant one.
$mergedDone = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "I also convert likert scale data*") {
        $p.Range.Text = "For 3 groups or more than 3 groups comparison, eg. in ethnicity/race group, I used ANOVA to conduct an analysis showing whether there is significant difference of extent of helpful among these categories."
        $nextP = $d.Paragraphs.Item($i + 1)
        $nextP.Range.Delete()
        $mergedDone = $true
        break
    }
}
if (-not $mergedDone) {
    Write-Output "NOT FOUND: I also convert likert scale data..."
}
